$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Point-sheet rows 2..16: renumber the "address" (F) column sequentially,
# collapse the "quality" (G) column to 1, and swap the demo "type"/"order"
# (H/I) values from FLOAT32/DCBA to INT16/AB.
For ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 6).Value = $r - 1
    $ws.Cells.Item($r, 7).Value = 1
    $ws.Cells.Item($r, 8).Value = "INT16"
    $ws.Cells.Item($r, 9).Value = "AB"
}

# Give the H column (type) its own explicit format (distinct xf/style)
# instead of sharing the plain body style.
$ws.Range("H2:H16").NumberFormat = "General"

# Move the active selection down to the row right after the table.
$ws.Range("H17").Select()
